# Automatische test-sync: 2025-06-30 19:41:50
# Adds the new "Testmail #3" log row (row 4) to the Logs sheet, adds the
# matching aggregate row to the Dashboard sheet, and extends the dashboard
# chart's category/value series references to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 4 with the new e-mail log entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Hoe kan ik een product retourneren?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #3: Hoe kan ik een product retourneren?"
$logs.Range("D4").Value = "Retour / Terugbetaling"
$logs.Range("E4").Value = "Beste klant,`nBedankt voor je vraag over het retourneren van een product. Om een product te retourneren, kun je contact opnemen met onze klantenservice via telefoon of e-mail. Zij zullen je voorzien van alle benodigde informatie en begeleiding betreffende het retourproces.`nMocht je verdere vragen hebben of als er meer informatie nodig is, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Je naam]`nE-mailassistent"
$logs.Range("F4").Value = "2025-06-30 19:41:48"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Nee"
$logs.Range("I4").Value = "Ja"
$logs.Range("J4").Value = "Nee"

# The multi-line content in E4 makes the engine auto-apply a custom row
# height; re-run AutoFit so the row goes back to the sheet's default
# height (no customHeight), matching the source log's plain row.
$logs.Rows.Item(4).EntireRow.AutoFit()

# Extend the existing conditional formatting blocks so they cover the new row
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "3")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "4")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count(); $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the matching aggregate row 4
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Retour / Terugbetaling"
$dash.Range("B4").Value = 1

# ---------------------------------------------------------------------
# 3. Update the chart's category/value series to span through row 4
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$4"
$series.Values = "='Dashboard'!`$B`$2:`$B`$4"
